$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Range("H28").Value = 1769.75
$ws.Range("I28").Value = 212.25
$ws.Range("J28").Value = 3327.25
$ws.Range("K28").Value = 212.25
$ws.Range("L28").Value = 3327.25
$ws.Range("M28").Value = 272.75
$ws.Range("N28").Value = -4297.25
# row 64
$ws.Range("I64").Value = 7878.75
$ws.Range("J64").Value = 8333.333000000001
$ws.Range("K64").Value = 7878.75
$ws.Range("L64").Value = 8333.333000000001
$ws.Range("M64").Value = -7630.75
$ws.Range("N64").Value = -8829.333000000001
# row 67
$ws.Range("I67").Value = 7878.75
$ws.Range("J67").Value = 8333.333000000001
$ws.Range("K67").Value = 7878.75
$ws.Range("L67").Value = 8333.333000000001
$ws.Range("M67").Value = -7020.75
$ws.Range("N67").Value = -10049.333
# row 69
$ws.Range("H69").Value = 7504.3335
$ws.Range("J69").Value = 10000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31748
# row 72
$ws.Range("H72").Value = 7504.3335
$ws.Range("J72").Value = 10000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98736
# row 86
$ws.Range("H86").Value = 5001.5
$ws.Range("I86").Value = 5001.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5001.5
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -3878.5
# row 89
$ws.Range("H89").Value = 5001.5
$ws.Range("I89").Value = 5001.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 25007.5
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -19391.5
# row 103
$ws.Range("H103").Value = 1375
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 1750
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 5250
$ws.Range("M103").Value = -2414
$ws.Range("N103").Value = -6422
# row 106
$ws.Range("H106").Value = 500000260
$ws.Range("I106").Value = 500000260
$ws.Range("K106").Value = 500000260
$ws.Range("M106").Value = -499999629
# row 138
$ws.Range("H138").Value = 4139.028
$ws.Range("J138").Value = 4378.037
$ws.Range("L138").Value = 13134.111
$ws.Range("N138").Value = -23414.111
# row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 432.66666
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 498
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 498
$ws.Range("M2").Value = -287
$ws.Range("N2").Value = -724
# row 35
$ws.Range("H35").Value = 2583.6
$ws.Range("I35").Value = 2583.6
$ws.Range("K35").Value = 2583.6
$ws.Range("M35").Value = -2177.6
# row 102
$ws.Range("H102").Value = 14001595
$ws.Range("I102").Value = 834743.7
$ws.Range("K102").Value = 834743.7
$ws.Range("M102").Value = -833121.7
# row 110
$ws.Range("H110").Value = 250000000
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
# row 116
$ws.Range("H116").Value = 432.66666
$ws.Range("I116").Value = 400
$ws.Range("J116").Value = 498
$ws.Range("K116").Value = 400
$ws.Range("L116").Value = 498
$ws.Range("M116").Value = 1894
$ws.Range("N116").Value = -5086

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 432.66666
$ws.Range("I3").Value = 400
$ws.Range("J3").Value = 498
$ws.Range("K3").Value = 400
$ws.Range("L3").Value = 498
$ws.Range("M3").Value = -286
$ws.Range("N3").Value = -726
# row 105
$ws.Range("H105").Value = 1933
$ws.Range("I105").Value = 1933
$ws.Range("K105").Value = 1933
$ws.Range("M105").Value = -186
# row 106
$ws.Range("H106").Value = 23902.166
$ws.Range("J106").Value = 23902.166
$ws.Range("L106").Value = 23902.166
$ws.Range("N106").Value = -26426.166
# row 134
$ws.Range("H134").Value = 2513.5715
$ws.Range("I134").Value = 1919
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5757
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3222
$ws.Range("N134").Value = -17070

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
# row 16
$ws.Range("H16").Value = 1098.8
$ws.Range("I16").Value = 799.5
$ws.Range("J16").Value = 1298.3334
$ws.Range("K16").Value = 799.5
$ws.Range("L16").Value = 1298.3334
$ws.Range("M16").Value = -512.5
$ws.Range("N16").Value = -1872.3334
# row 32
$ws.Range("H32").Value = 5933.5
$ws.Range("I32").Value = 3700
$ws.Range("K32").Value = 3700
$ws.Range("M32").Value = -3384
# row 35
$ws.Range("H35").Value = 944.5714
$ws.Range("I35").Value = 944.5714
$ws.Range("K35").Value = 944.5714
$ws.Range("M35").Value = -650.5714
# row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
# row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
# row 113
$ws.Range("H113").Value = 1098.8
$ws.Range("I113").Value = 799.5
$ws.Range("J113").Value = 1298.3334
$ws.Range("K113").Value = 799.5
$ws.Range("L113").Value = 1298.3334
$ws.Range("M113").Value = 1370.5
$ws.Range("N113").Value = -5638.3334
# row 132
$ws.Range("H132").Value = 1039
$ws.Range("I132").Value = 798.75
$ws.Range("K132").Value = 2396.25
$ws.Range("M132").Value = 133.75
# row 141
$ws.Range("H141").Value = 1488887.4
$ws.Range("J141").Value = 1488887.4
$ws.Range("L141").Value = 1488887.4
$ws.Range("N141").Value = -1499247.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 98
$ws.Range("H98").Value = 1525.4
$ws.Range("I98").Value = 423
$ws.Range("J98").Value = 2260.3333
$ws.Range("K98").Value = 1269
$ws.Range("L98").Value = 6780.999899999999
$ws.Range("M98").Value = 229
$ws.Range("N98").Value = -9776.999899999999
# row 107
$ws.Range("H107").Value = 247.5
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
# row 114
$ws.Range("H114").Value = 1564
$ws.Range("I114").Value = 1297.5
$ws.Range("J114").Value = 1640.1428
$ws.Range("K114").Value = 3892.5
$ws.Range("L114").Value = 4920.428400000001
$ws.Range("M114").Value = -638.5
$ws.Range("N114").Value = -11428.4284
# row 117
$ws.Range("H117").Value = 5654
$ws.Range("J117").Value = 7966.6665
$ws.Range("L117").Value = 23899.9995
$ws.Range("N117").Value = -30783.9995
# row 121
$ws.Range("H121").Value = 662.25
$ws.Range("I121").Value = 347.25
$ws.Range("J121").Value = 977.25
$ws.Range("K121").Value = 1041.75
$ws.Range("L121").Value = 2931.75
$ws.Range("M121").Value = 268.25
$ws.Range("N121").Value = -5551.75
# row 139
$ws.Range("H139").Value = 3703.7273
$ws.Range("I139").Value = 3624.1
$ws.Range("K139").Value = 10872.3
$ws.Range("M139").Value = -5732.299999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 250002260
$ws.Range("I70").Value = 4500
$ws.Range("K70").Value = 4500
$ws.Range("M70").Value = -4230
# row 73
$ws.Range("H73").Value = 250002260
$ws.Range("I73").Value = 4500
$ws.Range("K73").Value = 4500
$ws.Range("M73").Value = -3564
# row 80
$ws.Range("H80").Value = 19500
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("N80").Value = -3996
# row 83
$ws.Range("H83").Value = 19500
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 10000
$ws.Range("N83").Value = -19984
# row 102
$ws.Range("H102").Value = 1197.5
$ws.Range("I102").Value = 1197.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1197.5
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 424.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").ClearContents()
$ws.Range("N33").Value = 0
# row 68
$ws.Range("H68").Value = 8685.571
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 9999.833000000001
$ws.Range("K68").Value = 800
$ws.Range("L68").Value = 9999.833000000001
$ws.Range("M68").Value = -51
$ws.Range("N68").Value = -11497.833
# row 71
$ws.Range("H71").Value = 8685.571
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 9999.833000000001
$ws.Range("K71").Value = 4000
$ws.Range("L71").Value = 49999.165
$ws.Range("M71").Value = -256
$ws.Range("N71").Value = -57487.165
# row 122
$ws.Range("H122").Value = 4937.375
$ws.Range("I122").Value = 4749.8335
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 14249.5005
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -11799.5005
$ws.Range("N122").Value = -21400

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 750
$ws.Range("I132").Value = 750
$ws.Range("K132").Value = 2250
$ws.Range("M132").Value = 280
